$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.039.88"
$ws.Range("E2").Value = "  +5.35%  "
$ws.Range("D3").Value = "3.536.30"
$ws.Range("E3").Value = "  +5.27%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'188.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.07%  "
$ws.Range("D6").Value = "'562.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.86%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.08%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.526.72"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.631"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.98%  "
$ws.Range("D12").Value = "'54.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.82%  "
$ws.Range("D14").Value = "'9.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "4.105.71"
$ws.Range("E15").Value = "  +5.60%  "
$ws.Range("D16").Value = "3.547.33"
$ws.Range("E16").Value = "  +5.92%  "
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("D18").Value = "'18.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.60%  "
$ws.Range("D19").Value = "67.194.84"
$ws.Range("E19").Value = "  +5.62%  "
$ws.Range("D20").Value = "'12.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.10%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "'423.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.83%  "
$ws.Range("D23").Value = "'4.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.57%  "
$ws.Range("D24").Value = "'85.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").Value = "'4.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("D26").Value = "'11.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").Value = "'2.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.34%  "
$ws.Range("D28").Value = "'6.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "'12.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.61%  "
$ws.Range("D30").Value = "'8.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.54%  "
$ws.Range("D31").Value = "'30.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("D32").Value = "'631.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "'6.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").Value = "'11.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("D36").Value = "'60.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.94%  "
$ws.Range("D37").Value = "'38.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.99%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0810"
$ws.Range("E38").Value = "  +11.65%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.92%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'0.387"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").Value = "3.143.02"
$ws.Range("E42").Value = "  +6.86%  "
$ws.Range("D43").Value = "'3.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.11%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "'3.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.45%  "
$ws.Range("E47").Value = "  +9.79%  "
$ws.Range("D48").Value = "'0.0417"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.87%  "
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("E50").Value = "  +5.44%  "
$ws.Range("D51").Value = "'141.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.61%  "
